$wb = $excel.ActiveWorkbook
$srcWs = $wb.Worksheets.Item(1)
$sheetCount = $wb.Worksheets.Count
$lastWs = $wb.Worksheets.Item($sheetCount)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastWs)
$ws.Name = "cetificate"

# --- Header row ---
$ws.Range("B1").Value = 'name'
$ws.Range("C1").Value = 'cetificate'

# --- Column A: 0-based row index ---
for ($i = 0; $i -lt 110; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $i
}

# --- Column B: names, written top-to-bottom (drives shared-string order) ---
$names = @(
    'James Rolison', 'Jamie Chen', 'Jamie Conn', 'Jamie Hurley', 'Jan Macy-Buescher', 'Jan Manista', 'Jan Thede', 'Jane Greenfield', 'Janet Taffe', 'Janice Upton', 'Janine Dopson', 'Jasmine Geffner, CFA, CPA', 'Jason Crispin', 'Jason Crowley', 'Jason Ito', 'Jason Riley', 'Jason Rinne', 'Jason Schick', 'Jason Stegu', 'Jason Wells', 'Jay Chall', 'Jay Fort', 'Jay Massimo', 'JB Askew', 'JB Meanor II', 'JC Boyanton', 'JC Fanning', 'Jean Grasso', 'Jean-Paul (JP) Purdy, MBA', 'Jeanette Griffin', 'Jeanine Casey', 'Jeannette Lu, CFA', 'Jeff Bakalar', 'Jeff Billig', 'Jeff DeRosa', 'Jeff Erhardt', 'Jeff French', 'Jeff Maillet', 'Jeff Morrison', 'Jeff Norte', 'Jeff Ogden', 'Jeff Page', 'Jeff Rose', 'Jeff Skinner', 'Jeff Steele', 'Jeff Sullivan', 'Jeff Susman', 'Jeff Johnson', 'Jeffery Rose', 'Jeff Bryan', 'Jeffrey Clark', 'Jeffrey Hauser', 'Jeffrey Hoffman', 'Jeff Jacob', 'Jeffrey LaBauve', 'Jeffrey Leets', 'Jeffrey Miller', 'Jeffrey Mo', 'Jeffrey S. Ackerman', 'Jeffrey Saperstein', 'Jennifer Barrett', 'Jennifer Barrett', 'Jennifer Cann', 'Jennifer Choe', 'Jennifer Fitzgerald', 'Jennifer Fryhoff', 'Jennifer Heard', 'Jennifer Lyons', 'Jennifer Visconti', 'Jentri Smith', 'Jeremy Jackson', 'Jeremy Stump', 'Jeroen Fikke', 'Jerry Parisi, CFA', 'Jerry Wells', 'Jerry Zinkula, CFA', 'Jesper Lindquist', 'Jess Adams', 'Jesse Mason', 'Jessica Richardson', 'Jérôme Frizé', 'Jill Fedoruk', 'Jim Farner', 'Jimmy Simien', '金伟煌', 'Jo Ann Vasquez', 'Jo Ellen Bender', 'Joan Park, CFA', 'Joanne Nasuti', 'Jocelyn Boll', 'Joe Dancy', 'Joe Dougherty', 'Joe McCreery', 'Joe Netzel', 'Joel Makowsky', 'Joel Outlaw', 'Joey Powell', 'John Abate', 'John Brazzale', 'John Burda', 'John C. Thurston', 'John Canty', 'John Coffin', 'John Horst', 'John Lambert', 'John Tibe', 'John Wain', 'Joseph Onischuk', 'Joseph Sullivan', 'S. John Castellano'
)
for ($i = 0; $i -lt $names.Length; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $names[$i]
}

# --- Column C: cetificate, written top-to-bottom only where present ---
$certRows = @(5, 13, 19, 33, 41, 47, 75, 77, 89, 93)
$certVals = @('CPA', 'CPA, CFA', 'CPA', 'CFA', 'CPA', 'CFA', 'CFA', 'CFA', 'CFA', 'CPA, CFA')
for ($i = 0; $i -lt $certRows.Length; $i++) {
    $ws.Cells.Item($certRows[$i], 3).Value = $certVals[$i]
}

# --- Styling: reuse the existing bold/border/center style (style index 1) ---
$srcWs.Range("B1").Copy()
$ws.Range("B1:C1").PasteSpecial(-4122)
$srcWs.Range("A2").Copy()
$ws.Range("A2:A111").PasteSpecial(-4122)
$excel.CutCopyMode = 0
